$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (and thus the <sheet name=.../> element)
$ws.Name = "Through 2021-11-19"

# Update the label for the November row (shared string text)
$ws.Range("A12").Value = "November (through 11-19)"

# Update November row (row 12) values
$ws.Range("B12").Value = 19
$ws.Range("C12").Value = 43
$ws.Range("D12").Value = 78
$ws.Range("F12").Value = 30
$ws.Range("G12").Value = 117
$ws.Range("H12").Value = 131

# Update Total row (row 13) values
$ws.Range("B13").Value = 277
$ws.Range("C13").Value = 529
$ws.Range("D13").Value = 788
$ws.Range("F13").Value = 512
$ws.Range("G13").Value = 1174
$ws.Range("H13").Value = 1573
